$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''27.571.33'
$ws.Range('E2').Value = '  +5.85%  '
$ws.Range('D3').Value = '''1.724.03'
$ws.Range('E3').Value = '  +4.45%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '''225.93'
$ws.Range('E5').Value = '  +3.50%  '
$ws.Range('D6').Value = '''0.5365'
$ws.Range('E6').Value = '  +3.04%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '''0.2670'
$ws.Range('E8').Value = '  +1.26%  '
$ws.Range('D9').Value = '''0.06606'
$ws.Range('E9').Value = '  +4.45%  '
$ws.Range('E10').Value = '  +6.87%  '
$ws.Range('D11').Value = '''0.07723'
$ws.Range('E11').Value = '  +0.86%  '
$ws.Range('D12').Value = '''4.627'
$ws.Range('E12').Value = '  +0.96%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '''1.721.52'
$ws.Range('E13').Value = '  +3.71%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '''1.962.07'
$ws.Range('E14').Value = '  +4.51%  '
$ws.Range('D15').Value = '''0.5849'
$ws.Range('E15').Value = '  +4.56%  '
$ws.Range('D16').Value = '''0.0₅8311'
$ws.Range('E16').Value = '  +2.24%  '
$ws.Range('D17').Value = '''68.01'
$ws.Range('E17').Value = '  +4.31%  '
$ws.Range('D18').Value = '''27.576.55'
$ws.Range('E18').Value = '  +5.92%  '
$ws.Range('D19').Value = '''221.54'
$ws.Range('E19').Value = '  +15.60%  '
$ws.Range('D20').Value = '''1.004'
$ws.Range('E20').Value = '  +0.03%  '
$ws.Range('D21').Value = '''4.729'
$ws.Range('E21').Value = '  +2.37%  '
$ws.Range('E22').Value = '  +1.50%  '
$ws.Range('D23').Value = '''6.089'
$ws.Range('E23').Value = '  +3.23%  '
$ws.Range('E24').Value = '  +0.06%  '
$ws.Range('E25').Value = '  +2.93%  '
$ws.Range('D26').Value = '''1.738'
$ws.Range('E26').Value = '  +14.86%  '
$ws.Range('D27').Value = '''0.1235'
$ws.Range('E27').Value = '  +4.17%  '
$ws.Range('D28').Value = '''7.416'
$ws.Range('E28').Value = '  +3.13%  '
$ws.Range('D29').Value = '''16.62'
$ws.Range('E29').Value = '  +4.71%  '
$ws.Range('D30').Value = '''0.05582'
$ws.Range('E30').Value = '  +2.99%  '
$ws.Range('D31').Value = '''1.303'
$ws.Range('E31').Value = '  +2.45%  '
$ws.Range('D32').Value = '''3.572'
$ws.Range('E32').Value = '  +3.69%  '
$ws.Range('D33').Value = '''3.454'
$ws.Range('E33').Value = '  +3.34%  '
$ws.Range('D34').Value = '''1.661'
$ws.Range('E34').Value = '  +6.96%  '
$ws.Range('D35').Value = '''0.9702'
$ws.Range('E35').Value = '  +2.79%  '
$ws.Range('D36').Value = '''2.838'
$ws.Range('E36').Value = '  +2.06%  '
$ws.Range('D37').Value = '''2.427'
$ws.Range('E38').Value = '  +6.09%  '
$ws.Range('D39').Value = '''0.01652'
$ws.Range('D40').Value = '''5.919'
$ws.Range('E40').Value = '  +1.09%  '
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').Value = '''1.055.53'
$ws.Range('E41').Value = '  +2.98%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Value = '''0.8532'
$ws.Range('E42').Value = '  +3.55%  '
$ws.Range('D43').Value = '''1.004'
$ws.Range('D44').Value = '''101.39'
$ws.Range('E44').Value = '  +0.63%  '
$ws.Range('E45').Value = '  +4.62%  '
$ws.Range('D46').Value = '''0.0₈114'
$ws.Range('E46').Value = '  +2.87%  '
$ws.Range('D47').Value = '''59.19'
$ws.Range('E47').Value = '  +3.39%  '
$ws.Range('D48').Value = '''8.220'
$ws.Range('E48').Value = '  +3.81%  '
$ws.Range('D49').Value = '''0.4432'
$ws.Range('E49').Value = '  +2.30%  '
$ws.Range('E51').Value = '  +2.13%  '
